# Atualização automática via cronjob
# Refreshes the "vendas atipicas" report: the 7 existing data rows (2-8) get
# new values for the latest cronjob run, and 3 brand-new rows (9-11) are
# appended for newly detected atypical sales.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine text (no number/date auto-coercion) while
# keeping the cell's existing style untouched. Builds the text via a throwaway
# formula, then collapses it to a static value with Copy + PasteSpecial
# (xlPasteValues = -4163) so no NumberFormat bookkeeping leaks into styles.xml.
function Set-TextValue {
    param($Cell, [string]$Text)
    $escaped = $Text -replace '"', '""'
    $Cell.Formula = '="' + $escaped + '"'
    $Cell.Copy()
    $Cell.PasteSpecial(-4163)
}

# Helper: stamp a freshly-created cell with the same formatting as the
# existing index column (column A) cells above it (bold/centered/bordered),
# via Copy + PasteSpecial (xlPasteFormats = -4122).
function Copy-Format {
    param($FromCell, $ToCell)
    $FromCell.Copy()
    $ToCell.PasteSpecial(-4122)
}

# Column layout: A=index (no header, bordered/bold style), B=Dia, C=quantidade_atipica,
# D=cliente, E=id_produto, F=produto, G=estoque_atualizado, H=critico (boolean)
$rows = @(
    @{ R = 2;  A = 4;  Dia = "2025-04-17"; Qtd = 84;  Cliente = "METALURGICA SATO DA AMAZONIA LTDA";    IdProduto = "000032"; Produto = "LIMPADOR VEJA MULTIUSO GOLD 500ML";                       Estoque = 740;  Critico = $false },
    @{ R = 3;  A = 7;  Dia = "2025-04-22"; Qtd = 24;  Cliente = "BRAGA MOTOS LTDA";                     IdProduto = "000015"; Produto = "PANO MULTIUSO ROLO 28X300 M AZUL TALGE";                  Estoque = 0;    Critico = $false },
    @{ R = 4;  A = 8;  Dia = "2025-04-23"; Qtd = 400; Cliente = "V V REFEICOES LTDA";                   IdProduto = "000029"; Produto = "ESPONJA MULTIUSO JEITOSA";                                Estoque = 1121; Critico = $false },
    @{ R = 5;  A = 1;  Dia = "2025-04-24"; Qtd = 150; Cliente = "JURUA ESTALEIROS E NAVEGACAO LTDA";    IdProduto = "000088"; Produto = "VASSOURA PIACAVA 20 FUROS";                               Estoque = 2;    Critico = $false },
    @{ R = 6;  A = 3;  Dia = "2025-04-24"; Qtd = 300; Cliente = "MUSASHI DA AMAZONIA LTDA";             IdProduto = "000842"; Produto = "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND";               Estoque = 21;   Critico = $false },
    @{ R = 7;  A = 0;  Dia = "2025-04-28"; Qtd = 250; Cliente = "MAP SERVICOS DE CONSERVACAO - EIRELI"; IdProduto = "000098"; Produto = "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM";       Estoque = 683;  Critico = $false },
    @{ R = 8;  A = 11; Dia = "2025-04-28"; Qtd = 60;  Cliente = "MAP SERVICOS DE CONSERVACAO - EIRELI"; IdProduto = "000347"; Produto = "INSETICIDA BUZZOFF AEROSSOL 300ML";                       Estoque = 19;   Critico = $true  },
    @{ R = 9;  A = 5;  Dia = "2025-04-28"; Qtd = 70;  Cliente = "MAP SERVICOS DE CONSERVACAO - EIRELI"; IdProduto = "000349"; Produto = "DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO";   Estoque = 547;  Critico = $true  },
    @{ R = 10; A = 6;  Dia = "2025-04-28"; Qtd = 20;  Cliente = "LUCAS CLIENTE NOVO";                   IdProduto = "000158"; Produto = "AZULIM LIMPA CERAMICAS E AZULEJOS LAVANDA 5L 1:15 START"; Estoque = 9;    Critico = $true  },
    @{ R = 11; A = 9;  Dia = "2025-04-28"; Qtd = 250; Cliente = "MAP SERVICOS DE CONSERVACAO - EIRELI"; IdProduto = "000779"; Produto = "PEDRA SANITARIA NAFT PLUS FLORAL 25G";                    Estoque = -151; Critico = $false }
)

# New index-column cells (rows 9-11) need the same formatting (bold/border/
# center) as the existing rows above them before any values are written.
foreach ($row in $rows) {
    if ($row.R -gt 8) {
        Copy-Format $ws.Cells.Item(8, 1) $ws.Cells.Item($row.R, 1)
    }
}

# Write column-by-column (matching the source data pipeline's natural
# left-to-right, per-column regeneration order) so shared strings land in the
# same clustered order as the rest of the workbook.
foreach ($row in $rows) { $ws.Cells.Item($row.R, 1).Value = $row.A }
foreach ($row in $rows) { Set-TextValue $ws.Cells.Item($row.R, 2) $row.Dia }
foreach ($row in $rows) { $ws.Cells.Item($row.R, 3).Value = $row.Qtd }
foreach ($row in $rows) { Set-TextValue $ws.Cells.Item($row.R, 4) $row.Cliente }
foreach ($row in $rows) { Set-TextValue $ws.Cells.Item($row.R, 5) $row.IdProduto }
foreach ($row in $rows) { Set-TextValue $ws.Cells.Item($row.R, 6) $row.Produto }
foreach ($row in $rows) { $ws.Cells.Item($row.R, 7).Value = $row.Estoque }
foreach ($row in $rows) { $ws.Cells.Item($row.R, 8).Value = $row.Critico }
